# fix time_type database files
# Update quarterly income-statement figures on the "Overview" sheet
# (rows 11-37, columns D-H) with the corrected reported values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = -5196880
$ws.Range("E11").Value = 0
$ws.Range("G11").Value = 7965836
$ws.Range("H11").Value = 5768479
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("D13").Value = 4390595
$ws.Range("E13").Value = 4858810
$ws.Range("G13").Value = 8383926
$ws.Range("H13").Value = 5705797
$ws.Range("D14").Value = 5203777
$ws.Range("E14").Value = 0
$ws.Range("G14").Value = -8142767
$ws.Range("H14").Value = -5572552
$ws.Range("D15").Value = 4397492
$ws.Range("E15").Value = 4858810
$ws.Range("F15").Value = 4220330
$ws.Range("G15").Value = 3986665
$ws.Range("H15").Value = 5901724
$ws.Range("D16").Value = -2648839
$ws.Range("E16").Value = -2800829
$ws.Range("G16").Value = -5535697
$ws.Range("H16").Value = -4526724
$ws.Range("D17").Value = -187869
$ws.Range("E17").Value = -262965
$ws.Range("F17").Value = -233335
$ws.Range("G17").Value = -321735
$ws.Range("H17").Value = -271899
$ws.Range("D18").Value = 0
$ws.Range("E18").Value = 111251
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 12974
$ws.Range("H18").Value = 0
$ws.Range("D19").Value = -25242
$ws.Range("E19").Value = -30029
$ws.Range("F19").Value = -14099
$ws.Range("G19").Value = -21045
$ws.Range("H19").Value = -19396
$ws.Range("D20").Value = 1535542
$ws.Range("E20").Value = 1876238
$ws.Range("F20").Value = 1312457
$ws.Range("G20").Value = 781601
$ws.Range("H20").Value = 1083705
$ws.Range("D21").Value = 17191
$ws.Range("E21").Value = 21657
$ws.Range("G21").Value = 99848
$ws.Range("H21").Value = 32330
$ws.Range("D22").Value = 12722
$ws.Range("E22").Value = 57597
$ws.Range("F22").Value = 126000
$ws.Range("G22").Value = 123548
$ws.Range("H22").Value = 66358
$ws.Range("D24").Value = -89
$ws.Range("E24").Value = -864
$ws.Range("F24").Value = -557
$ws.Range("G24").Value = -78
$ws.Range("H24").Value = -75
$ws.Range("D25").Value = 1565366
$ws.Range("E25").Value = 1954628
$ws.Range("F25").Value = 1477460
$ws.Range("G25").Value = 965359
$ws.Range("H25").Value = 1182318
$ws.Range("D26").Value = -5510
$ws.Range("E26").Value = 352
$ws.Range("F26").Value = 0
$ws.Range("G26").Value = -11636
$ws.Range("H26").Value = -12302
$ws.Range("D27").Value = 1009
$ws.Range("E27").Value = -4361
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = -9583
$ws.Range("H27").Value = 9583
$ws.Range("D28").Value = 1560865
$ws.Range("E28").Value = 1950619
$ws.Range("F28").Value = 1477460
$ws.Range("G28").Value = 944140
$ws.Range("H28").Value = 1179599
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 0
$ws.Range("G29").Value = 0
$ws.Range("H29").Value = 0
$ws.Range("D30").Value = 1560865
$ws.Range("E30").Value = 1950619
$ws.Range("F30").Value = 1477460
$ws.Range("G30").Value = 944140
$ws.Range("H30").Value = 1179599
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("F35").Value = 1137
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 1300000
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0
$ws.Range("D37").Value = 1201
$ws.Range("E37").Value = 1500
$ws.Range("F37").Value = 1137
$ws.Range("G37").Value = 726
$ws.Range("H37").Value = 907